# Applies the "Knärot" appendix section (13 new paragraphs at the end of
# the body, before the final sectPr) and bumps the header date from
# 2023-09-13 to 2023-09-15.
$d = $word.ActiveDocument

# --- Phase 1: insert empty paragraph placeholders (correct style each) ---
# Doing this as a separate pass (before any text/italic formatting is set)
# avoids Word carrying the previous paragraph's trailing run formatting
# (e.g. italics) into the freshly inserted paragraph.
$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter()
$p0 = $d.Paragraphs($d.Paragraphs.Count)
$p0.Style = "Heading 1"
$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs($d.Paragraphs.Count)
$p1.Style = "Normal"
$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs($d.Paragraphs.Count)
$p2.Style = "Normal"
$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs($d.Paragraphs.Count)
$p3.Style = "Normal"
$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs($d.Paragraphs.Count)
$p4.Style = "Normal"
$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter()
$p5 = $d.Paragraphs($d.Paragraphs.Count)
$p5.Style = "Normal"
$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter()
$p6 = $d.Paragraphs($d.Paragraphs.Count)
$p6.Style = "Heading 2"
$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter()
$p7 = $d.Paragraphs($d.Paragraphs.Count)
$p7.Style = "Normal"
$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter()
$p8 = $d.Paragraphs($d.Paragraphs.Count)
$p8.Style = "Normal"
$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter()
$p9 = $d.Paragraphs($d.Paragraphs.Count)
$p9.Style = "Normal"
$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter()
$p10 = $d.Paragraphs($d.Paragraphs.Count)
$p10.Style = "Normal"
$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter()
$p11 = $d.Paragraphs($d.Paragraphs.Count)
$p11.Style = "Normal"
$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.InsertParagraphAfter()
$p12 = $d.Paragraphs($d.Paragraphs.Count)
$p12.Style = "Normal"

# --- Phase 2: fill in paragraph text, then italicise the quoted/citation runs ---
# Paragraph 0
$fullText = 'Knärot – ekologi samt krav på livsmiljön'
$p0.Range.Text = $fullText

# Paragraph 1
$fullText = 'Knärot är fridlyst enligt 8 och 15 §§ artskyddsförordningen och klassad som sårbar (VU) enligt rödlistan 2020. Knärot är beroende av hög och jämn luftfuktighet i gamla, ostörda skogsmiljöer och är känslig för snabba förändringar av ljus-/vindförhållanden eller uttorkning. På grund av ett alltför intensivt skogsbruk har den minskat med 40 (25-50) % under de senaste 60 åren och i framtiden bedöms minskningstakten uppgå till 30 (20-40) %. Till följd av att arten har en dokumenterat högre minskningstakt iförhållande till sin generationstid än vad som tidigare varit känt (data från Riksskogstaxeringen) höjdes den till hotkategori sårbar (VU) i rödlistan 2020 (Artdatabanken, 2021).'
$p1.Range.Text = $fullText

# Paragraph 2
$fullText = 'Samuel Johnsons doktorsavhandling “Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation“ (SLU, Uppsala 2014) visar att det krävs väl tilltagna skyddszoner för att knärotens växtplatser inte ska ta skada av skogsbruksåtgärder i intilliggande områden: “Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.” Vidare “More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).”'
$p2.Range.Text = $fullText
$pStart = $p2.Range.Start
$rng = $d.Range($pStart + 34, $pStart + 116)
$rng.Font.Italic = $true
$rng = $d.Range($pStart + 278, $pStart + 483)
$rng.Font.Italic = $true
$rng = $d.Range($pStart + 490, $pStart + 608)
$rng.Font.Italic = $true

# Paragraph 3
$fullText = 'Johnsons (2014) rekommendation på minst 50 meters breda skyddszoner runt knärotens växtplatser motsvarar en areal på 0,78 hektar, vilket ligger i linje med andra studier som gjorts på känsliga skogsarter: “In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).”'
$p3.Range.Text = $fullText
$pStart = $p3.Range.Start
$rng = $d.Range($pStart + 205, $pStart + 1070)
$rng.Font.Italic = $true

# Paragraph 4
$fullText = 'En nyligen publicerad vetenskaplig uppsats av Koelmeijer m.fl. (2022) inkluderar orkidén knärots skyddsbehov. I uppsatsen berörs problemet med uttorkning för växter, bl.a. för knärot, ett problem som blivit accentuerat på grund av den pågående klimatförändringen och torra somrar, t.ex. den exceptionellt torra sommaren 2018. I uppsatsen undersöks områden med tre olika avstånd från kalhyggeskant med avseende på skydd bl.a. för knärot. Det första området har avstånd upp till 20 m från hyggeskant (Strong edge effect), det andra 20 – 40 m från hyggeskant (Weak edge effect) och det tredje avser större avstånd från hyggeskant, där kanteffekten anses vara försumbar (Interior). Ett resultat var att man fann stor eller mycket stor uttorkningseffekt på känsliga och rödlistade skogsarter vid de kortare avstånden till hyggeskant, medan effekt av uttorkning inte konstaterades på större avstånd (Interior). För orkidén knärot fann man en rik förekomst (upp till 0,06 dm2/m2) på stort avstånd från hyggeskant (Interior), medan förekomsten var liten eller närmast försumbar i de områden som klassificerades som Weak edge effect respektive Strong edge effect. Arbetet påpekar att de allt oftare förekommande torra somrarna ger ytterligare skäl att utöka skyddsavståndet från hyggen till den fuktkrävande arten knärot (Koelmeijer m.fl., 2022).'
$p4.Range.Text = $fullText

# Paragraph 5
$fullText = 'Även Skogsstyrelsens egen vägledning för hänsyn till knärot ligger i linje med ovanstående forskningsstudier. Av vägledningen framgår det att för med hög sannolikhet kunna bevara befintliga förekomster krävs relativt stora avsättningar av uppvuxen skog med slutet och relativt tätt kronskikt. Som riktlinje kan krävas ett avstånd på 50 meter in från brynet för att vidmakthålla ett fungerande mikroklimat. Detta innebär att fristående hänsynsytor för många arter (kärlväxter, lavar och mossor) kan behöva ha en area överstigande 0,8 hektar (cirkelyta med radien 50 meter = 0,78 hektar) för att bibehålla lokalklimatet. Även ganska små förändringar i form av förändrade ljus- och fuktighetsförhållanden, till exempel till följd av gallring, kan leda till att arten försvinner till följd av konkurrens med mera ljuskrävande och snabbväxande arter (Skogsstyrelsen, 2022).'
$p5.Range.Text = $fullText

# Paragraph 6
$fullText = 'Referenser - knärot'
$p6.Range.Text = $fullText

# Paragraph 7
$fullText = 'de Graaf M & Roberts M.R., 2009. Short-term response of the herbaceous layer within leave patches after harvest. Forest Ecology and Management 257, 1014-1025'
$p7.Range.Text = $fullText
$pStart = $p7.Range.Start
$rng = $d.Range($pStart + 33, $pStart + 113)
$rng.Font.Italic = $true

# Paragraph 8
$fullText = 'Halpern, C. B., Halaj, J., Evans, S. A., & Dovciak, M., 2012. Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. Ecological Applications, 22, 2049-2064 '
$p8.Range.Text = $fullText
$pStart = $p8.Range.Start
$rng = $d.Range($pStart + 62, $pStart + 176)
$rng.Font.Italic = $true

# Paragraph 9
$fullText = 'Koelmeijer, I. A., Ehrlén, J., Jönsson, M., De Frenne, P., Berg, P., Andersson, J., Weibull, H. & Hylander, N. 2022. Interactive effects of drought and edge exposure on old-growth forest understory species. Landscape Ecology, 37, sid 1839-1853'
$p9.Range.Text = $fullText
$pStart = $p9.Range.Start
$rng = $d.Range($pStart + 117, $pStart + 207)
$rng.Font.Italic = $true

# Paragraph 10
$fullText = 'Rudolphi, J., Jönsson, M. T., & Gustafsson, L., 2014. Biological legacies buffer local species extinction after logging. Journal of Applied Ecology. 51, 53-62.'
$p10.Range.Text = $fullText
$pStart = $p10.Range.Start
$rng = $d.Range($pStart + 54, $pStart + 121)
$rng.Font.Italic = $true

# Paragraph 11
$fullText = 'Skogsstyrelsen, 2022. Vägledning för hänsyn till knärot. https://www.skogsstyrelsen.se/lag-och-tillsyn/artskydd/vagledningar-och-kunskapsstod-artskydd/vagledning-for-hansyn-till-knarot/'
$p11.Range.Text = $fullText
$pStart = $p11.Range.Start
$rng = $d.Range($pStart + 22, $pStart + 57)
$rng.Font.Italic = $true

# Paragraph 12
$fullText = 'SLU Artdatabanken, 2021. Artfaktablad. Naturvård – artfakta. SLU Artdatabanken, Uppsala '
$p12.Range.Text = $fullText
$pStart = $p12.Range.Start
$rng = $d.Range($pStart + 25, $pStart + 61)
$rng.Font.Italic = $true

# --- Update the date in the first-page header (2023-09-13 -> 2023-09-15) ---
$hdr = $d.Sections(1).Headers(2)  # wdHeaderFooterFirstPage
$hdrText = $hdr.Range.Text
$dateIdx = $hdrText.IndexOf("2023-09-13")
if ($dateIdx -ge 0) {
    $dateRng = $hdr.Range.Duplicate()
    $dateRng.Start = $dateIdx
    $dateRng.End = $dateIdx + 10
    $dateRng.Text = "2023-09-15"
}

